# Added RMA Test Cases: SO To inspection order SO to RMA Receipt
#
# The "RMA Details Maintenance Grid" sheet holds 3 sample RMA-receipt rows
# (row 2, 3, 4). Row 3's RMA / RMA-line / Id values are refreshed to a new
# generated test-case record (RMA-V2KN-*) while everything else on the row
# (Product, Qty, Reason for Return, Action, Receive To Site) stays as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Column E = Sales Order Line ("RMA#" value used for the order line lookup)
$ws.Range("E3").Value = "RMA-V2KN-001"

# Column F = Shipper Line
$ws.Range("F3").Value = "RMA-V2KN-1-1"

# Column J = Id (Salesforce record id of the RMA line)
$ws.Range("J3").Value = "a6h1K000000Q2JGQA0"
